# unitycar.xlsx - "fixed xls diff error"
#
# The Transmission sheet's brake-bias helper table (rows 63-73) used a plain
# linear formula in column F ( = E * $B$60 ) which could push the resulting
# G/H (1-F / 1+F) ratios outside the intended +/-60% clamp once E exceeded
# +/-1. The fix clamps F to $B$60 (keeping the sign of E) whenever E's
# magnitude would otherwise overshoot it.
#
# This also corrects a stray test value left in B65 (0.05 -> 0) and leaves
# the workbook with the Transmission sheet active/selected at C65, matching
# the state the author saved the file in.

$wb = $excel.ActiveWorkbook

$transmission = $wb.Worksheets.Item("Transmission")

# --- Fix the clamped brake-bias formula (column F, rows 63:73) -------------
# Old:  = E* $B$60
# New:  = IF(ABS(E*) > ABS($B$60), $B$60 * SIGN(E*), E*)
$transmission.Range("F63:F73").Formula = "=IF(ABS(E63) > ABS(`$B`$60), `$B`$60 * SIGN(E63), E63)"

# --- Correct the stray input value in B65 (0.05 -> 0) ----------------------
$transmission.Range("B65").Value = 0

# --- Restore the saved view state: Transmission tab active, C65 selected ---
$transmission.Activate()
$transmission.Range("C65").Select()
